# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (Overview!E2/F2, zh-cn!C2, de-de!C2)
#  - Latest Handback DateTime refreshed for zh-cn and de-de
#  - Error Detail (stale-handback warning) cleared now that files are in sync

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn (E2) and de-de (F2) status columns
$wsOverview.Cells.Item(2, 5).Value = $newStatus
$wsOverview.Cells.Item(2, 6).Value = $newStatus

# zh-cn detail sheet
$wsZhCn.Cells.Item(2, 3).Value = $newStatus                     # Status
$wsZhCn.Cells.Item(2, 11).Value = "2016-08-28 02:49:09"         # Latest Handback DateTime
$wsZhCn.Cells.Item(2, 16).Value = ""                            # Error Detail

# de-de detail sheet
$wsDeDe.Cells.Item(2, 3).Value = $newStatus                     # Status
$wsDeDe.Cells.Item(2, 11).Value = "2016-08-28 02:49:16"         # Latest Handback DateTime
$wsDeDe.Cells.Item(2, 16).Value = ""                            # Error Detail
